$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data per diff
$updates = @(
    ,@('D2', '42.508.38')
    ,@('E2', '  +1.25%  ')
    ,@('D3', '2.303.97')
    ,@('E3', '  +0.41%  ')
    ,@('E4', '  +0.17%  ')
    ,@('D5', '316.56')
    ,@('E5', '  -0.26%  ')
    ,@('D6', '103.83')
    ,@('E6', '  +0.17%  ')
    ,@('E7', '  -0.17%  ')
    ,@('E8', '  +0.22%  ')
    ,@('E9', '  +0.20%  ')
    ,@('D10', '40.05')
    ,@('E10', '  +1.72%  ')
    ,@('E11', '  +0.07%  ')
    ,@('D12', '8.53')
    ,@('E12', '  +3.15%  ')
    ,@('E13', '  +0.97%  ')
    ,@('D14', '0.997')
    ,@('E14', '  +3.91%  ')
    ,@('D15', '15.33')
    ,@('E15', '  +0.65%  ')
    ,@('D16', '2.653.90')
    ,@('E16', '  +0.71%  ')
    ,@('D17', '2.298.11')
    ,@('E17', '  +0.33%  ')
    ,@('D18', '42.599.06')
    ,@('E18', '  +1.58%  ')
    ,@('D19', '7.61')
    ,@('E19', '  +3.17%  ')
    ,@('E20', '  +0.23%  ')
    ,@('D21', '13.65')
    ,@('E21', '  +34.42%  ')
    ,@('D22', '73.97')
    ,@('E22', '  +0.94%  ')
    ,@('E23', '  -2.27%  ')
    ,@('D24', '267.18')
    ,@('E24', '  -3.83%  ')
    ,@('E25', '  -1.00%  ')
    ,@('D26', '1.01')
    ,@('E26', '  -0.19%  ')
    ,@('D27', '10.92')
    ,@('E27', '  +1.29%  ')
    ,@('D28', '2.34')
    ,@('E28', '  -0.99%  ')
    ,@('D29', '22.60')
    ,@('E29', '  -0.94%  ')
    ,@('D30', '38.17')
    ,@('E30', '  +6.25%  ')
    ,@('D31', '6.57')
    ,@('E31', '  +13.21%  ')
    ,@('D32', '165.63')
    ,@('E32', '  +1.55%  ')
    ,@('D33', '0.0880')
    ,@('E33', '  +1.12%  ')
    ,@('D35', '2.64')
    ,@('E35', '  -7.01%  ')
    ,@('E36', '  +0.26%  ')
    ,@('E37', '  +2.29%  ')
    ,@('D38', '0.0353')
    ,@('E38', '  +2.07%  ')
    ,@('E39', '  -0.25%  ')
    ,@('D40', '2.77')
    ,@('E40', '  -1.80%  ')
    ,@('E41', '  +13.71%  ')
    ,@('D42', '98.09')
    ,@('E42', '  -1.18%  ')
    ,@('D43', '70.06')
    ,@('E43', '  +1.03%  ')
    ,@('E44', '  +0.45%  ')
    ,@('E45', '  +0.29%  ')
    ,@('B46', 'Celestia')
    ,@('C46', 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia')
    ,@('D46', '12.35')
    ,@('E46', '  +3.65%  ')
    ,@('B47', 'Aave')
    ,@('C47', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave')
    ,@('D47', '116.89')
    ,@('E47', '  +3.71%  ')
    ,@('D48', '80.43')
    ,@('E48', '  +4.28%  ')
    ,@('D49', '1.644.68')
    ,@('E49', '  +4.04%  ')
    ,@('E50', '  +0.50%  ')
    ,@('E51', '  +0.12%  ')
)

foreach ($u in $updates) {
    $cellRef = $u[0]
    $newVal = $u[1]
    $col = $cellRef.Substring(0, 1)
    $range = $ws.Range($cellRef)
    if ($col -eq "D" -or $col -eq "E") {
        # Force text number format so numeric-looking / percent-looking strings
        # are preserved exactly as text (matching original inlineStr cells)
        $range.NumberFormat = "@"
    }
    $range.Value = $newVal
}
